$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying CI run (with time_filter + first_treshold_filter now wired up) produced
# refreshed counts for the long-format Volatility/Trend pivot summary. Re-write the affected
# rows (year columns + total_count_of_occurrences + percentage_of_occurrences) for each symbol
# block: EURUSD (rows 4-9), GBPUSD (rows 13-18), USDJPY (rows 22-27), XAUUSD (rows 31-36).

function Set-RowValues($row, $startCol, $values) {
    $col = $startCol
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# ---- EURUSD block (years 2016-2022 in B:H, total in I, pct in J) ----
# Row 4: High Volatility + No Trend
Set-RowValues 4 2 @(52, 54, 90, 42, 89, 59, 96, 482, 26.41095890410959)
# Row 6: High Volatility + Trend False
Set-RowValues 6 2 @(13, 26, 15, 23, 10, 9, 37, 133, 7.287671232876712)
# Row 7: Low Volatility + No Trend
Set-RowValues 7 2 @(104, 76, 80, 105, 58, 104, 35, 562, 30.7945205479452)
# Row 9: Low Volatility + Trend False
Set-RowValues 9 2 @(23, 37, 12, 57, 11, 22, 9, 171, 9.36986301369863)

# ---- GBPUSD block (years 2017-2022 in B:G, total in H, pct in I) ----
# Row 13: High Volatility + No Trend
Set-RowValues 13 2 @(45, 69, 57, 82, 42, 97, 392, 25.04792332268371)
# Row 15: High Volatility + Trend False
Set-RowValues 15 2 @(32, 25, 4, 31, 14, 28, 134, 8.562300319488818)
# Row 16: Low Volatility + No Trend
Set-RowValues 16 2 @(109, 63, 99, 59, 115, 46, 491, 31.37380191693291)
# Row 18: Low Volatility + Trend False
Set-RowValues 18 2 @(23, 26, 2, 25, 49, 8, 133, 8.498402555910543)

# ---- USDJPY block (years 2017-2022 in B:G, total in H, pct in I) ----
# Row 22: High Volatility + No Trend
Set-RowValues 22 2 @(52, 58, 68, 66, 75, 82, 401, 25.62300319488818)
# Row 23: High Volatility + Trend (only the 2019 count, total and pct moved)
$ws.Cells.Item(23, 4).Value = 8
$ws.Cells.Item(23, 8).Value = 263
$ws.Cells.Item(23, 9).Value = 16.80511182108626
# Row 24: High Volatility + Trend False
Set-RowValues 24 2 @(9, 8, 32, 10, 28, 14, 101, 6.453674121405751)
# Row 25: Low Volatility + No Trend
Set-RowValues 25 2 @(105, 106, 102, 65, 80, 28, 486, 31.05431309904154)
# Row 26: Low Volatility + Trend (only the 2019 count, total and pct moved)
$ws.Cells.Item(26, 4).Value = 3
$ws.Cells.Item(26, 8).Value = 177
$ws.Cells.Item(26, 9).Value = 11.30990415335463
# Row 27: Low Volatility + Trend False
Set-RowValues 27 2 @(28, 15, 46, 18, 24, 6, 137, 8.753993610223642)

# ---- XAUUSD block (years 2016-2022 in B:H, total in I, pct in J); column C (2017) untouched ----
# Row 31: High Volatility + No Trend
$ws.Cells.Item(31, 2).Value = 62
$ws.Cells.Item(31, 4).Value = 63
$ws.Cells.Item(31, 5).Value = 73
$ws.Cells.Item(31, 6).Value = 76
$ws.Cells.Item(31, 7).Value = 59
$ws.Cells.Item(31, 8).Value = 88
$ws.Cells.Item(31, 9).Value = 489
$ws.Cells.Item(31, 10).Value = 27.07641196013289
# Row 33: High Volatility + Trend False
$ws.Cells.Item(33, 2).Value = 24
$ws.Cells.Item(33, 4).Value = 31
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = 27
$ws.Cells.Item(33, 7).Value = 10
$ws.Cells.Item(33, 8).Value = 8
$ws.Cells.Item(33, 9).Value = 114
$ws.Cells.Item(33, 10).Value = 6.312292358803987
# Row 34: Low Volatility + No Trend
$ws.Cells.Item(34, 2).Value = 77
$ws.Cells.Item(34, 4).Value = 68
$ws.Cells.Item(34, 5).Value = 75
$ws.Cells.Item(34, 6).Value = 61
$ws.Cells.Item(34, 7).Value = 96
$ws.Cells.Item(34, 9).Value = 533
$ws.Cells.Item(34, 10).Value = 29.51273532668882
# Row 36: Low Volatility + Trend False
$ws.Cells.Item(36, 2).Value = 17
$ws.Cells.Item(36, 4).Value = 39
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 16
$ws.Cells.Item(36, 7).Value = 10
$ws.Cells.Item(36, 9).Value = 108
$ws.Cells.Item(36, 10).Value = 5.980066445182724
